# Update NATMI LR-pair metrics (Bdnf-Sort1) with newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 0.4664623486936776
$ws.Range("J2").Value = 0.4664623486936776
$ws.Range("M2").Value = 0.6574793333333333
$ws.Range("N2").Value = 1.972438
$ws.Range("O2").Value = 0.04234443143670402
$ws.Range("P2").Value = 0.04234443143670403
$ws.Range("Q2").Value = 0.2222205632342222
$ws.Range("R2").Value = 1.999985069108
$ws.Range("S2").Value = 0.01975208294206336
$ws.Range("T2").Value = 0.01975208294206336
$ws.Range("I3").Value = 0.4664623486936776
$ws.Range("J3").Value = 0.4664623486936776
$ws.Range("O3").Value = 0.1192373589365509
$ws.Range("P3").Value = 0.119237358936551
$ws.Range("Q3").Value = 0.6257491755686666
$ws.Range("R3").Value = 5.631742580117999
$ws.Range("S3").Value = 0.05561973850157462
$ws.Range("T3").Value = 0.05561973850157463
$ws.Range("I4").Value = 0.4664623486936776
$ws.Range("J4").Value = 0.4664623486936776
$ws.Range("M4").Value = 5.370269333333333
$ws.Range("N4").Value = 16.110808
$ws.Range("O4").Value = 0.3458679080132824
$ws.Range("P4").Value = 0.3458679080132824
$ws.Range("Q4").Value = 1.815090171614222
$ws.Range("R4").Value = 16.335811544528
$ws.Range("S4").Value = 0.1613343567096445
$ws.Range("T4").Value = 0.1613343567096445
$ws.Range("I5").Value = 0.4664623486936776
$ws.Range("J5").Value = 0.4664623486936776
$ws.Range("M5").Value = 1.801189666666667
$ws.Range("N5").Value = 5.403569
$ws.Range("O5").Value = 0.1160041821512257
$ws.Range("P5").Value = 0.1160041821512257
$ws.Range("Q5").Value = 0.6087816938504444
$ws.Range("R5").Value = 5.479035244654
$ws.Range("S5").Value = 0.05411158326454994
$ws.Range("T5").Value = 0.05411158326454994
$ws.Range("I6").Value = 0.4664623486936776
$ws.Range("J6").Value = 0.4664623486936776
$ws.Range("M6").Value = 5.846608
$ws.Range("N6").Value = 17.539824
$ws.Range("O6").Value = 0.3765461194622369
$ws.Range("P6").Value = 0.376546119462237
$ws.Range("Q6").Value = 1.976087242442667
$ws.Range("R6").Value = 17.784785181984
$ws.Range("S6").Value = 0.1756445872758451
$ws.Range("T6").Value = 0.1756445872758452
$ws.Range("G7").Value = 0.38659
$ws.Range("H7").Value = 1.15977
$ws.Range("I7").Value = 0.5335376513063224
$ws.Range("J7").Value = 0.5335376513063224
$ws.Range("M7").Value = 0.6574793333333333
$ws.Range("N7").Value = 1.972438
$ws.Range("O7").Value = 0.04234443143670402
$ws.Range("P7").Value = 0.04234443143670403
$ws.Range("Q7").Value = 0.2541749354733333
$ws.Range("R7").Value = 2.28757441926
$ws.Range("S7").Value = 0.02259234849464066
$ws.Range("T7").Value = 0.02259234849464067
$ws.Range("G8").Value = 0.38659
$ws.Range("H8").Value = 1.15977
$ws.Range("I8").Value = 0.5335376513063224
$ws.Range("J8").Value = 0.5335376513063224
$ws.Range("O8").Value = 0.1192373589365509
$ws.Range("P8").Value = 0.119237358936551
$ws.Range("Q8").Value = 0.7157292466899999
$ws.Range("R8").Value = 6.441563220209999
$ws.Range("S8").Value = 0.06361762043497632
$ws.Range("T8").Value = 0.06361762043497633
$ws.Range("G9").Value = 0.38659
$ws.Range("H9").Value = 1.15977
$ws.Range("I9").Value = 0.5335376513063224
$ws.Range("J9").Value = 0.5335376513063224
$ws.Range("M9").Value = 5.370269333333333
$ws.Range("N9").Value = 16.110808
$ws.Range("O9").Value = 0.3458679080132824
$ws.Range("P9").Value = 0.3458679080132824
$ws.Range("Q9").Value = 2.076092421573333
$ws.Range("R9").Value = 18.68483179416
$ws.Range("S9").Value = 0.1845335513036378
$ws.Range("T9").Value = 0.1845335513036378
$ws.Range("G10").Value = 0.38659
$ws.Range("H10").Value = 1.15977
$ws.Range("I10").Value = 0.5335376513063224
$ws.Range("J10").Value = 0.5335376513063224
$ws.Range("M10").Value = 1.801189666666667
$ws.Range("N10").Value = 5.403569
$ws.Range("O10").Value = 0.1160041821512257
$ws.Range("P10").Value = 0.1160041821512257
$ws.Range("Q10").Value = 0.6963219132366667
$ws.Range("R10").Value = 6.26689721913
$ws.Range("S10").Value = 0.06189259888667577
$ws.Range("T10").Value = 0.06189259888667577
$ws.Range("G11").Value = 0.38659
$ws.Range("H11").Value = 1.15977
$ws.Range("I11").Value = 0.5335376513063224
$ws.Range("J11").Value = 0.5335376513063224
$ws.Range("M11").Value = 5.846608
$ws.Range("N11").Value = 17.539824
$ws.Range("O11").Value = 0.3765461194622369
$ws.Range("P11").Value = 0.376546119462237
$ws.Range("Q11").Value = 2.26024018672
$ws.Range("R11").Value = 20.34216168048
$ws.Range("S11").Value = 0.2009015321863918
$ws.Range("T11").Value = 0.2009015321863918
